$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 192 (existing row): update changed values and fill newly reported columns -----
$ws.Range("B192").Value = 34869.4
$ws.Range("C192").Value = 14970
$ws.Range("D192").Value = 7063.4
$ws.Range("E192").Value = 30240.1
$ws.Range("F192").Value = 6650.9
$ws.Range("G192").Value = 15573.9
$ws.Range("H192").Value = 3133.6
$ws.Range("I192").Value = 4877.4
$ws.Range("J192").Value = 1533.1
$ws.Range("K192").Value = 70027.10000000001
$ws.Range("L192").Value = 4069.6
$ws.Range("M192").Value = 1620
$ws.Range("N192").Value = 17313.8
$ws.Range("O192").Value = 1391.7
$ws.Range("P192").Value = 75687.5
$ws.Range("Q192").Value = 113583
$ws.Range("R192").Value = 1345.4
$ws.Range("S192").Value = 51598.5
$ws.Range("T192").Value = 18653.3

# ----- New rows 193-197 (new trading days) -----
# The "Serie" column holds a dd-mm-yyyy-looking label that must stay plain TEXT.
# Force text format before assignment (and clear the temporary format again
# afterwards) so Excel does not silently convert it into a date serial number.

# Row 193
$ws.Range("A193").NumberFormat = "@"
$ws.Range("A193").Value = "28-09-2021"
$ws.Range("A193").ClearFormats()
$ws.Range("B193").Value = 34300
$ws.Range("C193").Value = 14546.7
$ws.Range("D193").Value = 7028.1
$ws.Range("E193").Value = 30184
$ws.Range("F193").Value = 6506.5
$ws.Range("G193").Value = 15248.6
$ws.Range("H193").Value = 3097.9
$ws.Range("I193").Value = 4883.8
$ws.Range("J193").Value = 1546.8
$ws.Range("K193").Value = 69440.5
$ws.Range("L193").Value = 4109.7
$ws.Range("M193").Value = 1616.5
$ws.Range("N193").Value = 17181.4
$ws.Range("O193").Value = 1383.8
$ws.Range("P193").Value = 75808.39999999999
$ws.Range("Q193").Value = 110123.9
$ws.Range("R193").Value = 1350.5
$ws.Range("S193").Value = 50925.5
$ws.Range("T193").Value = 18515.3

# Row 194
$ws.Range("A194").NumberFormat = "@"
$ws.Range("A194").Value = "29-09-2021"
$ws.Range("A194").ClearFormats()
$ws.Range("B194").Value = 34390.7
$ws.Range("C194").Value = 14512.4
$ws.Range("D194").Value = 7108.2
$ws.Range("E194").Value = 29544.3
$ws.Range("F194").Value = 6560.8
$ws.Range("G194").Value = 15365.3
$ws.Range("H194").Value = 3060.3
$ws.Range("I194").Value = 4833.9
$ws.Range("J194").Value = 1547.7
$ws.Range("K194").Value = 69871.39999999999
$ws.Range("L194").Value = 4045.2
$ws.Range("M194").Value = 1617
$ws.Range("N194").Value = 16855.5
$ws.Range("O194").Value = 1391.9
$ws.Range("P194").Value = 76589.3
$ws.Range("Q194").Value = 111106.8
$ws.Range("R194").Value = 1360.4
$ws.Range("S194").Value = 51084.6
$ws.Range("T194").Value = 18458.6

# Row 195
$ws.Range("A195").NumberFormat = "@"
$ws.Range("A195").Value = "30-09-2021"
$ws.Range("A195").ClearFormats()
$ws.Range("B195").Value = 33843.9
$ws.Range("C195").Value = 14448.6
$ws.Range("D195").Value = 7086.4
$ws.Range("E195").Value = 29452.7
$ws.Range("F195").Value = 6520
$ws.Range("G195").Value = 15260.7
$ws.Range("H195").Value = 3068.8
$ws.Range("I195").Value = 4866.4
$ws.Range("J195").Value = 1537.8
$ws.Range("K195").Value = 70340.89999999999
$ws.Range("L195").Value = 4079.5
$ws.Range("M195").Value = 1605.7
$ws.Range("N195").Value = 16934.8
$ws.Range("O195").Value = 1406.4
$ws.Range("P195").Value = 77363.60000000001
$ws.Range("Q195").Value = 110979.1
$ws.Range("R195").Value = 1362
$ws.Range("S195").Value = 51385.6
$ws.Range("T195").Value = 18279.2

# Row 196
$ws.Range("A196").NumberFormat = "@"
$ws.Range("A196").Value = "01-10-2021"
$ws.Range("A196").ClearFormats()
$ws.Range("B196").Value = 34326.5
$ws.Range("C196").Value = 14566.7
$ws.Range("D196").Value = 7027.1
$ws.Range("E196").Value = 28771.1
$ws.Range("F196").Value = 6517.7
$ws.Range("G196").Value = 15156.4
$ws.Range("H196").Value = 3019.2
$ws.Range("J196").Value = 1524.5
$ws.Range("K196").Value = 70812
$ws.Range("L196").Value = 4078.6
$ws.Range("M196").Value = 1605.2
$ws.Range("N196").Value = 16570.9
$ws.Range("O196").Value = 1401.5
$ws.Range("P196").Value = 77626.2
$ws.Range("Q196").Value = 112899.6
$ws.Range("R196").Value = 1370.3
$ws.Range("S196").Value = 51060.1
$ws.Range("T196").Value = 18348.7

# Row 197
$ws.Range("A197").NumberFormat = "@"
$ws.Range("A197").Value = "04-10-2021"
$ws.Range("A197").ClearFormats()
$ws.Range("D197").Value = 7022.8
$ws.Range("E197").Value = 28444.9
$ws.Range("F197").Value = 6509.6
$ws.Range("G197").Value = 15119.6
$ws.Range("J197").Value = 1522.5
$ws.Range("K197").Value = 71105.3
$ws.Range("L197").Value = 4124.6
$ws.Range("M197").Value = 1614.5
$ws.Range("N197").Value = 16408.4
$ws.Range("O197").Value = 1404.7
